# Auto-applied scheduled market-data refresh for Levequest profit tables.
# Updates currentAveragePrice* and derived Leve Profit columns (H-N) per sheet,
# mirroring a fresh Universalis market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 40556.855
$ws.Range("H70").Value = 11824.75
$ws.Range("J70").Value = 15366.333
$ws.Range("L70").Value = 46098.999
$ws.Range("N70").Value = -46638.999
$ws.Range("H72").Value = 40556.855
$ws.Range("H73").Value = 11824.75
$ws.Range("J73").Value = 15366.333
$ws.Range("L73").Value = 46098.999
$ws.Range("N73").Value = -47970.999
$ws.Range("H82").Value = 1626.5555
$ws.Range("I82").Value = 1626.5555
$ws.Range("K82").Value = 4879.666499999999
$ws.Range("M82").Value = -4473.666499999999
$ws.Range("H85").Value = 1626.5555
$ws.Range("I85").Value = 1626.5555
$ws.Range("K85").Value = 4879.666499999999
$ws.Range("M85").Value = -3475.666499999999
$ws.Range("H97").Value = 642.5
$ws.Range("J97").Value = 642.5
$ws.Range("L97").Value = 1927.5
$ws.Range("N97").Value = -2919.5
$ws.Range("H99").Value = 539
$ws.Range("J99").Value = 485
$ws.Range("L99").Value = 1455
$ws.Range("N99").Value = -4451
$ws.Range("H100").Value = 5612.154
$ws.Range("I100").Value = 4295.3335
$ws.Range("J100").Value = 6740.857
$ws.Range("K100").Value = 4295.3335
$ws.Range("L100").Value = 6740.857
$ws.Range("M100").Value = -3754.3335
$ws.Range("N100").Value = -7822.857
$ws.Range("H107").Value = 1368.1
$ws.Range("J107").Value = 401.44446
$ws.Range("L107").Value = 401.44446
$ws.Range("N107").Value = -4241.44446
$ws.Range("H141").Value = 5274.852
$ws.Range("I141").Value = 4383.727
$ws.Range("J141").Value = 9195.799999999999
$ws.Range("K141").Value = 13151.181
$ws.Range("L141").Value = 27587.4
$ws.Range("M141").Value = -7971.181
$ws.Range("N141").Value = -37947.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2333.3333
$ws.Range("H74").Value = 1347.95
$ws.Range("I74").Value = 1361.8387
$ws.Range("K74").Value = 1361.8387
$ws.Range("M74").Value = -487.8387
$ws.Range("H77").Value = 1347.95
$ws.Range("I77").Value = 1361.8387
$ws.Range("K77").Value = 6809.1935
$ws.Range("M77").Value = -2441.1935
$ws.Range("H102").Value = 2887.4092
$ws.Range("I102").Value = 2896.15
$ws.Range("K102").Value = 2896.15
$ws.Range("M102").Value = -1274.15
$ws.Range("H132").Value = 3166
$ws.Range("I132").Value = 3166
$ws.Range("K132").Value = 9498
$ws.Range("M132").Value = -6968

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 500
$ws.Range("I24").Value = 500
$ws.Range("K24").Value = 500
$ws.Range("M24").Value = -265
$ws.Range("H134").Value = 3093.3635
$ws.Range("I134").Value = 3032.5715
$ws.Range("K134").Value = 9097.7145
$ws.Range("M134").Value = -6562.7145
$ws.Range("H138").Value = 69999.586
$ws.Range("J138").Value = 69999.586
$ws.Range("L138").Value = 69999.586
$ws.Range("N138").Value = -80279.586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1707.5555
$ws.Range("J16").Value = 1699.8
$ws.Range("L16").Value = 1699.8
$ws.Range("N16").Value = -2273.8
$ws.Range("H31").Value = 1947.125
$ws.Range("I31").Value = 1925.2858
$ws.Range("K31").Value = 1925.2858
$ws.Range("M31").Value = -1630.2858
$ws.Range("H34").Value = 1947.125
$ws.Range("I34").Value = 1925.2858
$ws.Range("K34").Value = 1925.2858
$ws.Range("M34").Value = -1723.2858
$ws.Range("H35").Value = 503
$ws.Range("I35").Value = 453.75
$ws.Range("J35").Value = 700
$ws.Range("K35").Value = 453.75
$ws.Range("L35").Value = 700
$ws.Range("M35").Value = -159.75
$ws.Range("N35").Value = -1288
$ws.Range("H39").Value = 7799.6
$ws.Range("I39").Value = 7799.6
$ws.Range("K39").Value = 7799.6
$ws.Range("M39").Value = -7408.6
$ws.Range("H49").Value = 7799.6
$ws.Range("I49").Value = 7799.6
$ws.Range("K49").Value = 7799.6
$ws.Range("M49").Value = -7617.6
$ws.Range("H105").Value = 761
$ws.Range("I105").Value = 725.125
$ws.Range("J105").Value = 856.6667
$ws.Range("K105").Value = 725.125
$ws.Range("L105").Value = 856.6667
$ws.Range("M105").Value = 1021.875
$ws.Range("N105").Value = -4350.6667
$ws.Range("H113").Value = 1707.5555
$ws.Range("J113").Value = 1699.8
$ws.Range("L113").Value = 1699.8
$ws.Range("N113").Value = -6039.8
$ws.Range("H132").Value = 3214.6667
$ws.Range("I132").Value = 2478.6155
$ws.Range("J132").Value = 7999
$ws.Range("K132").Value = 7435.8465
$ws.Range("L132").Value = 23997
$ws.Range("M132").Value = -4905.8465
$ws.Range("N132").Value = -29057
$ws.Range("H134").Value = 8373.25
$ws.Range("I134").Value = 8373.25
$ws.Range("K134").Value = 25119.75
$ws.Range("M134").Value = -22584.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 4975
$ws.Range("J95").Value = 4975
$ws.Range("L95").Value = 14925
$ws.Range("N95").Value = -19043
$ws.Range("H122").Value = 7940.4
$ws.Range("I122").Value = 880.8
$ws.Range("K122").Value = 7927.2
$ws.Range("M122").Value = -5477.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 23332.666
$ws.Range("I33").Value = 18000
$ws.Range("J33").Value = 25999
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 25999
$ws.Range("M33").Value = -17748
$ws.Range("N33").Value = -26503
$ws.Range("H38").Value = 18495.5
$ws.Range("J38").Value = 18495.5
$ws.Range("L38").Value = 18495.5
$ws.Range("N38").Value = -19421.5
$ws.Range("H40").Value = 500
$ws.Range("J40").Value = 500
$ws.Range("L40").Value = 500
$ws.Range("N40").Value = -802
$ws.Range("H47").Value = 21666.334
$ws.Range("J47").Value = 21666.334
$ws.Range("L47").Value = 21666.334
$ws.Range("N47").Value = -22802.334
$ws.Range("H122").Value = 4152.9287
$ws.Range("I122").Value = 1073.8572
$ws.Range("J122").Value = 7232
$ws.Range("K122").Value = 3221.5716
$ws.Range("L122").Value = 21696
$ws.Range("M122").Value = -771.5715999999998
$ws.Range("N122").Value = -26596
$ws.Range("H132").Value = 8323.429
$ws.Range("I132").Value = 7419
$ws.Range("K132").Value = 22257
$ws.Range("M132").Value = -19727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 998
$ws.Range("J55").Value = 997.5
$ws.Range("L55").Value = 997.5
$ws.Range("N55").Value = -1343.5
$ws.Range("H132").Value = 2856.5264
$ws.Range("I132").Value = 1948.1428
$ws.Range("K132").Value = 5844.428400000001
$ws.Range("M132").Value = -3314.428400000001
$ws.Range("H135").Value = 72194.62
$ws.Range("J135").Value = 72194.62
$ws.Range("L135").Value = 72194.62
$ws.Range("N135").Value = -82334.62

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3473.5833
$ws.Range("I132").Value = 2435.375
$ws.Range("K132").Value = 7306.125
$ws.Range("M132").Value = -4776.125
